# Update workbook per commit "Update data: 2025-11-12 15:30"
# 1. Remove the "distance from Dma50" worksheet entirely.
# 2. Bump the "Last Updated" timestamp on the Metadata sheet.
# 3. Refresh Price / % Change figures on the Stock List sheet.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- 1. Delete the "distance from Dma50" sheet ---
$wb.Worksheets("distance from Dma50").Delete()

# --- 2. Update the Metadata "Last Updated" timestamp ---
$meta = $wb.Worksheets("Metadata")
$meta.Range("A2").Value = "12 Nov 2025, 03:28 PM"

# --- 3. Update Stock List Price (D) and % Change (E) values ---
$ws = $wb.Worksheets("Stock List")

$ws.Range("D2").Value = 3021.2
$ws.Range("E2").Value = 2.0124
$ws.Range("D3").Value = 123.1
$ws.Range("E3").Value = -0.5413
$ws.Range("D4").Value = 423.75
$ws.Range("E4").Value = 2.0224
$ws.Range("D5").Value = 549.1
$ws.Range("E5").Value = -0.1546
$ws.Range("D6").Value = 1435.5
$ws.Range("E6").Value = -1.4147
$ws.Range("D7").Value = 665.95
$ws.Range("E7").Value = -3.0146
$ws.Range("D8").Value = 9.86
$ws.Range("E8").Value = 0.4073
$ws.Range("D9").Value = 1095.1
$ws.Range("E9").Value = -0.9855
$ws.Range("D16").Value = 626.9
$ws.Range("E16").Value = 0.9663
$ws.Range("D17").Value = 313.4
$ws.Range("E17").Value = -0.4289
$ws.Range("D18").Value = 1651
$ws.Range("E18").Value = -0.8706
$ws.Range("D21").Value = 606.2
$ws.Range("E21").Value = -2.6185
$ws.Range("D23").Value = 142.6
$ws.Range("E23").Value = 5.7158
$ws.Range("D25").Value = 72.51
$ws.Range("E25").Value = 1.1015
$ws.Range("D26").Value = 229.29
$ws.Range("E26").Value = -3.5137
$ws.Range("D27").Value = 218.14
$ws.Range("E27").Value = 3.6886
$ws.Range("D28").Value = 423.65
$ws.Range("E28").Value = -2.6987
$ws.Range("D29").Value = 310.65
$ws.Range("E29").Value = -2.2806
$ws.Range("D30").Value = 65.63
$ws.Range("E30").Value = -1.3379
$ws.Range("D31").Value = 346.2
$ws.Range("E31").Value = -2.4789
$ws.Range("D32").Value = 626.9
$ws.Range("E32").Value = -2.0928
$ws.Range("D33").Value = 290
$ws.Range("E33").Value = -4.0371
$ws.Range("D34").Value = 651.55
$ws.Range("E34").Value = 5.6682
$ws.Range("D35").Value = 284.15
$ws.Range("E35").Value = 0.7446
$ws.Range("D36").Value = 1010.2
$ws.Range("E36").Value = 3.8019
$ws.Range("D38").Value = 467.7
$ws.Range("E38").Value = -0.5528
$ws.Range("D39").Value = 322.5
$ws.Range("E39").Value = 8.659
$ws.Range("D40").Value = 67.89
$ws.Range("E40").Value = 2.167
$ws.Range("D41").Value = 342
$ws.Range("E41").Value = 2.2269
$ws.Range("D42").Value = 218.4
$ws.Range("E42").Value = -0.2922
$ws.Range("D43").Value = 43.53
$ws.Range("E43").Value = -0.4118
$ws.Range("D44").Value = 142.2
$ws.Range("E44").Value = -2.6761
$ws.Range("D46").Value = 11.6
$ws.Range("E46").Value = -0.3436
$ws.Range("D51").Value = 23.67
$ws.Range("E51").Value = 1.5444
$ws.Range("D52").Value = 105.3
$ws.Range("E52").Value = 2.6016
$ws.Range("D53").Value = 154.36
$ws.Range("E53").Value = 0.7309
$ws.Range("D55").Value = 70.02
$ws.Range("E55").Value = 0.7482
$ws.Range("D56").Value = 289.5
$ws.Range("E56").Value = 6.7084
$ws.Range("D57").Value = 205.5
$ws.Range("E57").Value = 1.773
$ws.Range("D58").Value = 222
$ws.Range("E58").Value = -0.2785
$ws.Range("D59").Value = 316.15
$ws.Range("E59").Value = 0
$ws.Range("D61").Value = 89.22
$ws.Range("E61").Value = 2.7525
$ws.Range("D62").Value = 591.05
$ws.Range("E62").Value = -2.0224
$ws.Range("D65").Value = 122.19
$ws.Range("E65").Value = -0.2856
$ws.Range("D69").Value = 170.78
$ws.Range("E69").Value = -0.7093
$ws.Range("D70").Value = 436.55
$ws.Range("E70").Value = 0.1261
$ws.Range("D72").Value = 1496.8
$ws.Range("E72").Value = -0.2133
$ws.Range("D73").Value = 228.55
$ws.Range("E73").Value = -0.665
$ws.Range("D74").Value = 123.1
$ws.Range("E74").Value = -0.3239
$ws.Range("D75").Value = 80
$ws.Range("E75").Value = -0.8183
$ws.Range("D76").Value = 221.5
$ws.Range("E76").Value = 2.4846
